# Listas sem duplicação de professores
# Remove duplicated teacher list entries, leaving a simple "-" placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "-"
$ws.Range("B19").Value = "-"
$ws.Range("F19").Value = "-"
$ws.Range("B20").Value = "-"
